# Update the TraceabilityMatrix worksheet with newly-traced test plan /
# software design references that were filled in for several requirements.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TraceabilityMatrix")

# Row 10 (requirement 3.1.1.1.2.2): link to software design AccountController.
$ws.Range("F10").Value = "3.1.1.3"
$ws.Range("G10").Value = "AccountController"

# Row 11 (requirement 3.1.1.1.2.3): new test plan entry.
$ws.Range("C11").Value = "5.2.1.20"
$ws.Range("D11").Value = "Graph View Tests"

# Row 20 (requirement 3.1.1.1.4.3): new test plan entry.
$ws.Range("C20").Value = "5.2.1.21"
$ws.Range("D20").Value = "Patient/Physician Association Tests"

# Row 21 (requirement 3.1.1.1.4.4): link to existing test plan entry.
$ws.Range("C21").Value = "5.2.1.9"
$ws.Range("D21").Value = "Username/Password Retrieval Tests"

# Row 22 (requirement 3.1.1.1.4.5): link to existing test plan entry.
$ws.Range("C22").Value = "5.2.1.4"
$ws.Range("D22").Value = "Account editing Tests"

# Row 35 (requirement 3.1.1.7): link to existing test plan entry.
$ws.Range("C35").Value = "5.2.1.5"
$ws.Range("D35").Value = "Account Management Tests"

# Row 37 (requirement 3.1.1.7.2): link to existing test plan entry.
$ws.Range("C37").Value = "5.2.1.5"
$ws.Range("D37").Value = "Account Management Tests"

# Row 48 (requirement 3.1.4.1): link to existing test plan entry. The test
# plan title cell picks up the same style already used for the rest of
# column D, so copy formats from a neighbouring, already-styled cell first.
$ws.Range("D12").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("C48").Value = "5.2.1.14"
$ws.Range("D48").Value = "Experiment Tests"

# Row 49 (requirement 4.1.4.1.1): link to existing test plan entry.
$ws.Range("D12").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("C49").Value = "5.2.1.14"
$ws.Range("D49").Value = "Experiment Tests"

$excel.CutCopyMode = 0

# Move the active selection to where the editor left off.
$ws.Range("C21").Select()
